# Generate Report for Handoff
# Updates the "Priority" column (E) and the handoff/handback timestamp
# columns for the rows whose handoff/handback has just completed.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 10, 12, 13, 14)

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-09-02 00:29:51"
}

# --- zh-cn sheet: Priority (column E) + Latest Handback DateTime (column H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
    $wsZhCn.Cells.Item($r, 8).Value = "2016-09-02 00:29:46"
}

# --- de-de sheet: Priority (column E) + Latest Handback DateTime (column H) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
    $wsDeDe.Cells.Item($r, 8).Value = "2016-09-02 00:29:51"
}
